$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "29.187.48"
Set-TextCell "D3" "1.827.04"
Set-TextCell "E4" "  +0.27%  "
Set-TextCell "D5" "233.72"
Set-TextCell "E5" "  -2.34%  "
Set-TextCell "E6" "  -4.50%  "
Set-TextCell "E7" "  +0.25%  "
Set-TextCell "D8" "0.06960"
Set-TextCell "E8" "  -5.85%  "
Set-TextCell "D9" "0.2749"
Set-TextCell "E9" "  -4.95%  "
Set-TextCell "D10" "23.25"
Set-TextCell "E10" "  -6.46%  "
Set-TextCell "D11" "0.07602"
Set-TextCell "E11" "  -1.48%  "
Set-TextCell "D12" "1.834.51"
Set-TextCell "E12" "  -0.18%  "
Set-TextCell "D13" "4.753"
Set-TextCell "E13" "  -4.29%  "
Set-TextCell "D14" "0.6257"
Set-TextCell "E14" "  -6.75%  "
Set-TextCell "D15" "0.000009663"
Set-TextCell "E15" "  -5.99%  "
Set-TextCell "D16" "78.22"
Set-TextCell "E16" "  -4.34%  "
Set-TextCell "D17" "28.835.11"
Set-TextCell "E17" "  -1.87%  "
Set-TextCell "D18" "5.704"
Set-TextCell "E18" "  -9.01%  "
Set-TextCell "D19" "220.71"
Set-TextCell "E19" "  -5.98%  "
Set-TextCell "E20" "  +0.22%  "
Set-TextCell "D21" "11.54"
Set-TextCell "E21" "  -6.14%  "
Set-TextCell "D22" "6.857"
Set-TextCell "E22" "  -5.84%  "
Set-TextCell "E23" "  -0.15%  "
Set-TextCell "D24" "155.40"
Set-TextCell "E24" "  -1.09%  "
Set-TextCell "D25" "7.953"
Set-TextCell "E25" "  -6.22%  "
Set-TextCell "D26" "0.1288"
Set-TextCell "E26" "  -4.37%  "
Set-TextCell "D27" "16.51"
Set-TextCell "E27" "  -4.81%  "
Set-TextCell "D28" "0.06512"
Set-TextCell "E28" "  -10.77%  "
Set-TextCell "D29" "1.450"
Set-TextCell "E29" "  -2.84%  "
Set-TextCell "D30" "1.437"
Set-TextCell "E30" "  -2.58%  "
Set-TextCell "D31" "3.836"
Set-TextCell "E31" "  -4.80%  "
Set-TextCell "D32" "3.752"
Set-TextCell "E32" "  -7.11%  "
Set-TextCell "D33" "1.092"
Set-TextCell "E33" "  -5.80%  "
Set-TextCell "D34" "1.718"
Set-TextCell "E34" "  -5.42%  "
Set-TextCell "D35" "0.6451"
Set-TextCell "E35" "  -8.97%  "
Set-TextCell "D36" "2.540"
Set-TextCell "E36" "  -1.59%  "
Set-TextCell "D37" "2.730"
Set-TextCell "E37" "  -2.01%  "
Set-TextCell "D38" "0.01740"
Set-TextCell "E38" "  -5.28%  "
Set-TextCell "D39" "6.522"
Set-TextCell "E39" "  -3.91%  "
Set-TextCell "D40" "1.169.76"
Set-TextCell "E40" "  -5.15%  "
Set-TextCell "D41" "0.8922"
Set-TextCell "E41" "  -6.48%  "
Set-TextCell "D42" "1.003"
Set-TextCell "E42" "  +0.21%  "
Set-TextCell "D43" "1.979.30"
Set-TextCell "E43" "  -0.61%  "
Set-TextCell "D44" "100.39"
Set-TextCell "E44" "  -0.72%  "
Set-TextCell "D45" "62.02"
Set-TextCell "E45" "  -5.04%  "
Set-TextCell "E46" "  -2.57%  "
Set-TextCell "D47" "1.585"
Set-TextCell "E47" "  -6.72%  "
Set-TextCell "B48" "Cronos"
Set-TextCell "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D48" "0.05545"
Set-TextCell "E48" "  -2.04%  "
Set-TextCell "B49" "EnergySwap"
Set-TextCell "C49" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D49" "8.428"
Set-TextCell "E49" "  -5.07%  "
Set-TextCell "D50" "0.4546"
Set-TextCell "E50" "  -0.64%  "
Set-TextCell "D51" "0.3640"
Set-TextCell "E51" "  -6.23%  "
